# Revert "Modificar Bankia por Caixabank": restore "Bankia" wherever the
# document currently reads "Grupo CaixaBank".
$d = $word.ActiveDocument

# 1) Global replace: "Grupo CaixaBank" -> "Bankia" (7 occurrences across the doc).
#    wdReplaceAll = 2, MatchWildcards not needed.
$d.Content.Find.Execute("Grupo CaixaBank", $false, $false, $false, $false, $false, $true, 1, $false, "Bankia", 2)

# 2) The "ALTA EN NEO CLIENTES" paragraph also regains the extra clause about
#    Bankia Habitat, which the forward ("Bankia -> Caixabank") commit had
#    dropped. Replace the now-shorter sentence with the restored full text.
$d.Content.Find.Execute(
    "La Oficina, deberá dar de alta como cliente Bankia al comprador/es de la operación si la oferta aprobada corresponde un activo/lote cuyo propietario sea Bankia.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "La Oficina, deberá dar de alta como cliente Bankia al comprador/es de la operación si la oferta aprobada corresponde un activo/lote cuyo propietario sea Bankia, o bien, en caso de corresponder a propietario distinto de Bankia, darlo/s de alta como cliente Bankia Habitat.",
    2)
